$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for 6 new product rows (the list gained 6 new drugs and is
#        kept sorted alphabetically, so everything from the old row 9 onward
#        cascades down). Insert 6 blank rows right after the current last
#        data row (16), which pushes the totals/footer block from 17-18 to
#        23-24. ---
$ws.Rows("17:22").Insert()

# --- 2. The newly inserted rows have no formatting yet; clone it from the
#        last existing data row (16) so the new rows look identical
#        (borders/fonts/number formats/merges match the rest of the table). ---
$ws.Range("A16:M16").Copy()
$ws.Range("A17:M22").PasteSpecial(-4122)
$ws.Range("N16").Copy()
$ws.Range("N17:N22").PasteSpecial(-4122)

# --- 3. Row heights: rows 17-22 follow the same alternating auto-height
#        pattern as the rest of the table, and the totals/footer rows need
#        their (slightly different) final heights. ---
$ws.Rows("17").RowHeight = 25.5
$ws.Rows("18").RowHeight = 25.5
$ws.Rows("19").RowHeight = 24.75
$ws.Rows("20").RowHeight = 25.5
$ws.Rows("21").RowHeight = 24.75
$ws.Rows("22").RowHeight = 25.5
$ws.Rows("23").RowHeight = 25.5
$ws.Rows("24").RowHeight = 17.25

# --- 4. Re-create the 3 merged groups (name / balance / price) for each of
#        the 6 new rows, same layout as every other data row. ---
for ($r = 17; $r -le 22; $r++) {
    $ws.Range("B$r`:G$r").Merge()
    $ws.Range("H$r`:K$r").Merge()
    $ws.Range("L$r`:M$r").Merge()
}

# --- 5. Rewrite the whole product table (rows 4-22) with the new,
#        re-sorted data set -- 6 new medicines were added and the rows
#        re-sequenced alphabetically by name. ---
$data = @(
    @(1,  "ALVEOLIN-P SYRUP 100 ML", "1:0", 50, 1),
    @(2,  "ANTINAL 220MG/5ML 60ML SUSP.", "1:0", 48, 2),
    @(3,  "APEXIDONE 4MG 30 F.C.TAB.", "0:2", 42, 0.33),
    @(4,  "ATROVENT 250MCG/2ML 20 UNIT DOSE VIAL", "1:19", 42.9, 0.15),
    @(5,  "DEPAKINE CHRONO 500MG 30 SCORED PROLONGED REL. F.C. TAB.", "1:0", 139.68, 1),
    @(6,  "DOLIPRANE 1 GM 15 TABS.", "10:1", 32, 0.67),
    @(7,  "GOURYST 0.5 MG 100 TABS.", "0:6", 19, 0.1),
    @(8,  "MEGAFEN-N 100MG/5ML SUSP. 120 ML", "1:0", 35, 1),
    @(9,  "MUCOPHYLLINE SYRUP 125 ML", "10:0", 50, 1),
    @(10, "NEUROVIT 6 I.M. AMPS", "3:4", 11, 0.17),
    @(11, "NEXIUM 20MG 28 F.C. TAB.", "0:0", 332, 1),
    @(12, "PANADOL ADVANCE 500 MG 48 TABLETS", "2:3", 23, 0.25),
    @(13, "PULMICORT 0.25MG/ML 20 NEBULIZER VIAL SUSP.", "0:19", 169.2, 0.3),
    @(14, "SPASMO-DIGESTIN 30 TABS.", "3:1", 73.32, 1),
    @(15, "TELFAST 120MG 20 F.C. TAB", "0:0", 116, 1),
    @(16, "TUSSISTOP 60 MG 20 TABS.", "0:0", 30, 0.5),
    @(17, "VOLTAREN 75MG/3ML 3 AMP.", "4:1", 17, 0.33),
    @(18, "YEAST MEPACO 60 TABS", "1:0", 60, 1),
    @(19, "سرنجات 3 سم", "-2:0", 4, 2)
)

$r = 4
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 8).Value = $row[2]
    $ws.Cells.Item($r, 12).Value = $row[3]
    $ws.Cells.Item($r, 14).Value = $row[4]
    $r = $r + 1
}

# --- 6. Update the grand total of the "سعر البيع" column. ---
$ws.Range("K23").Value = 1294.1
